# Apply updated odds for Jogos_da_Semana_FlashScore_2025-05-23.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 1.62
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 15
$ws.Range("T2").Value = 17
$ws.Range("U2").Value = 29
$ws.Range("Y2").Value = 41
$ws.Range("AA2").Value = 8
$ws.Range("AB2").Value = 15
$ws.Range("AF2").Value = 8.5

# Row 5
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 1.04
$ws.Range("K5").Value = 13
$ws.Range("V5").Value = 8.5
$ws.Range("AI5").Value = 26

# Row 8
$ws.Range("G8").Value = 2.7
$ws.Range("I8").Value = 2.45
$ws.Range("N8").Value = 1.9
$ws.Range("O8").Value = 1.9
$ws.Range("T8").Value = 9
$ws.Range("U8").Value = 13
$ws.Range("AB8").Value = 13

# Row 13
$ws.Range("J13").Value = 1.03
$ws.Range("L13").Value = 1.19

# Row 14
$ws.Range("G14").Value = 1.91
$ws.Range("H14").Value = 3.45
$ws.Range("L14").Value = 1.28
$ws.Range("M14").Value = 3.05
$ws.Range("N14").Value = 1.82
$ws.Range("O14").Value = 1.78
$ws.Range("P14").Value = 1.38
$ws.Range("Q14").Value = 2.57
$ws.Range("R14").Value = 1.72
$ws.Range("S14").Value = 1.9
$ws.Range("T14").Value = 7.2
$ws.Range("U14").Value = 9
$ws.Range("V14").Value = 8.5
$ws.Range("X14").Value = 15.5
$ws.Range("Y14").Value = 27
$ws.Range("Z14").Value = 10.25
$ws.Range("AA14").Value = 6.7
$ws.Range("AH14").Value = 50
$ws.Range("AJ14").Value = 37

# Row 15
$ws.Range("G15").Value = 3.75
$ws.Range("H15").Value = 3.6
$ws.Range("I15").Value = 1.9
$ws.Range("J15").Value = 1.03
$ws.Range("K15").Value = 12
$ws.Range("L15").Value = 1.19
$ws.Range("M15").Value = 4
$ws.Range("N15").Value = 1.75
$ws.Range("O15").Value = 2.05
$ws.Range("AB15").Value = 13
$ws.Range("AD15").Value = 151
$ws.Range("AE15").Value = 8.5
$ws.Range("AF15").Value = 10
$ws.Range("AH15").Value = 17

# Row 18
$ws.Range("J18").Value = 1.02
$ws.Range("K18").Value = 15
$ws.Range("L18").Value = 1.17
$ws.Range("N18").Value = 1.67
$ws.Range("O18").Value = 2.15

# Row 23
$ws.Range("R23").Value = 1.77
$ws.Range("S23").Value = 1.87

# Row 26
$ws.Range("J26").Value = 1.03
$ws.Range("L26").Value = 1.22
$ws.Range("R26").Value = 1.7

# Row 27
$ws.Range("G27").Value = 1.42
$ws.Range("H27").Value = 5
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 21
$ws.Range("AA27").Value = 10
$ws.Range("AB27").Value = 17
$ws.Range("AD27").Value = 151
$ws.Range("AJ27").Value = 41

# Row 37
$ws.Range("L37").Value = 1.25
$ws.Range("M37").Value = 3.75
$ws.Range("N37").Value = 1.8
$ws.Range("O37").Value = 2
$ws.Range("P37").Value = 1.33

# Row 38
$ws.Range("R38").Value = 1.92
$ws.Range("S38").Value = 1.77

# Row 40
$ws.Range("G40").Value = 1.95
$ws.Range("H40").Value = 3.5
$ws.Range("I40").Value = 3.4
$ws.Range("O40").Value = 1.83
$ws.Range("W40").Value = 17
$ws.Range("X40").Value = 15
$ws.Range("Y40").Value = 25
$ws.Range("Z40").Value = 10.75
$ws.Range("AA40").Value = 6.8
$ws.Range("AB40").Value = 14
$ws.Range("AF40").Value = 18.5
$ws.Range("AG40").Value = 11.75
$ws.Range("AI40").Value = 30
$ws.Range("AJ40").Value = 37

# Row 41
$ws.Range("G41").Value = 1.3
$ws.Range("H41").Value = 5
$ws.Range("I41").Value = 8.5
$ws.Range("L41").Value = 1.2
$ws.Range("M41").Value = 3.65
$ws.Range("N41").Value = 1.6
$ws.Range("O41").Value = 2.07
$ws.Range("R41").Value = 2
$ws.Range("S41").Value = 1.65
$ws.Range("T41").Value = 7.1
$ws.Range("U41").Value = 6.2
$ws.Range("W41").Value = 7.7
$ws.Range("Z41").Value = 13
$ws.Range("AA41").Value = 10.25
$ws.Range("AD41").Value = 1000
$ws.Range("AE41").Value = 22
$ws.Range("AF41").Value = 60
$ws.Range("AH41").Value = 250

# Row 42
$ws.Range("L42").Value = 1.3
$ws.Range("N42").Value = 1.88
$ws.Range("O42").Value = 1.72
$ws.Range("P42").Value = 1.42
$ws.Range("Q42").Value = 2.47
$ws.Range("T42").Value = 9.25
$ws.Range("U42").Value = 15
$ws.Range("AD42").Value = 500

# Row 43
$ws.Range("G43").Value = 1.65
$ws.Range("I43").Value = 4.7
$ws.Range("T43").Value = 7.3
$ws.Range("U43").Value = 8
$ws.Range("X43").Value = 12.5
$ws.Range("Y43").Value = 24
$ws.Range("AB43").Value = 15.5
$ws.Range("AC43").Value = 70
$ws.Range("AE43").Value = 13.5
$ws.Range("AF43").Value = 28
$ws.Range("AG43").Value = 15.5

# Row 44
$ws.Range("H44").Value = 2.87
$ws.Range("I44").Value = 2.4
$ws.Range("T44").Value = 6.1
$ws.Range("AC44").Value = 70

# Row 45
$ws.Range("H45").Value = 3.65
$ws.Range("AB45").Value = 15.5
$ws.Range("AC45").Value = 75
$ws.Range("AE45").Value = 11.25

# Row 48
$ws.Range("P48").Value = 1.36
$ws.Range("R48").Value = 1.69

# Row 49
$ws.Range("J49").Value = 1.05
$ws.Range("K49").Value = 11
$ws.Range("N49").Value = 1.85
$ws.Range("O49").Value = 1.95

# Row 50
$ws.Range("G50").Value = 1.7
$ws.Range("H50").Value = 3.25
$ws.Range("I50").Value = 4.5
$ws.Range("K50").Value = 8.5
$ws.Range("N50").Value = 2.1
$ws.Range("O50").Value = 1.7
$ws.Range("U50").Value = 7.5
$ws.Range("W50").Value = 13
$ws.Range("X50").Value = 15
$ws.Range("Y50").Value = 29
$ws.Range("Z50").Value = 8.5
$ws.Range("AF50").Value = 23
$ws.Range("AG50").Value = 17

# Row 51
$ws.Range("G51").Value = 3.2
$ws.Range("H51").Value = 3.2
$ws.Range("Q51").Value = 2.57
$ws.Range("T51").Value = 9
$ws.Range("U51").Value = 16.5
$ws.Range("W51").Value = 40
$ws.Range("X51").Value = 29
$ws.Range("Y51").Value = 40
$ws.Range("Z51").Value = 8.75
$ws.Range("AA51").Value = 6.2
$ws.Range("AB51").Value = 15
$ws.Range("AC51").Value = 75
$ws.Range("AE51").Value = 7
$ws.Range("AG51").Value = 9
$ws.Range("AH51").Value = 21

# Row 53
$ws.Range("G53").Value = 1.9
$ws.Range("H53").Value = 3.5
$ws.Range("I53").Value = 3.45
$ws.Range("M53").Value = 3.25
$ws.Range("N53").Value = 1.87
$ws.Range("O53").Value = 1.83
$ws.Range("Q53").Value = 2.72
$ws.Range("T53").Value = 7.3
$ws.Range("AD53").Value = 600
$ws.Range("AE53").Value = 10.5
$ws.Range("AF53").Value = 18.5
$ws.Range("AG53").Value = 12

# Row 54
$ws.Range("G54").Value = 1.6
$ws.Range("H54").Value = 4.1
$ws.Range("I54").Value = 5
$ws.Range("L54").Value = 1.24
$ws.Range("M54").Value = 3.8
$ws.Range("N54").Value = 1.72
$ws.Range("O54").Value = 2.05
$ws.Range("P54").Value = 1.35
$ws.Range("Q54").Value = 3
$ws.Range("W54").Value = 12.5
$ws.Range("X54").Value = 13
$ws.Range("Y54").Value = 26
$ws.Range("AE54").Value = 13.5
$ws.Range("AF54").Value = 32
$ws.Range("AG54").Value = 17
$ws.Range("AH54").Value = 100

# Row 57
$ws.Range("P57").Value = 1.44
$ws.Range("Q57").Value = 2.63

# Row 68
$ws.Range("N68").Value = 1.5
$ws.Range("O68").Value = 2.5
